$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 to push the old data down by one row,
# keeping the original (now moved) configuration data intact.
$ws.Rows.Item(2).Insert()

# Row 1 becomes the "x"/"y" column headers
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"

# Row 2 becomes the "default" configuration row
$ws.Range("A2").Value = "default"
$ws.Range("B2").Value = "default"

# Move the active selection to C1
$ws.Range("C1").Select()
